$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rows 22 and 23 swap places (same match date, two different fixtures
#        that were listed in the wrong order). Columns A-E (index, pais,
#        torneio, temporada, data_partida) stay put; F..V (match detail
#        columns) swap between the two rows. ---
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
foreach ($col in $cols) {
    $v22 = $ws.Range($col + "22").Value()
    $v23 = $ws.Range($col + "23").Value()
    $ws.Range($col + "22").Value = $v23
    $ws.Range($col + "23").Value = $v22
}

# --- 2) Append a new match row 44 (index 43), copying the formatting from
#        the last existing row (43) so styles (bold/border index column,
#        date-time number format column) line up. ---
$ws.Range("A43:V43").Copy()
$ws.Range("A44").PasteSpecial(-4122)

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "lebanon"
$ws.Range("C44").Value = "premier-league"
$ws.Range("D44").Value = "2023-2024"
$ws.Range("E44").Value = 45240.55208333334
$ws.Range("F44").Value = "Al Ansar"
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = "Al Ghazieh"
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 1.17
$ws.Range("K44").Value = "09/11/2023 01:42"
$ws.Range("L44").Value = 1.16
$ws.Range("M44").Value = "10/11/2023 12:16"
$ws.Range("N44").Value = 6.19
$ws.Range("O44").Value = "09/11/2023 01:42"
$ws.Range("P44").Value = 7.1
$ws.Range("Q44").Value = "10/11/2023 12:16"
$ws.Range("R44").Value = 10.71
$ws.Range("S44").Value = "09/11/2023 01:42"
$ws.Range("T44").Value = 14.24
$ws.Range("U44").Value = "10/11/2023 12:16"
$ws.Range("V44").Value = "https://www.betexplorer.com/football/lebanon/premier-league/al-ansar-al-ghazieh/KxNLGPGH/"
